$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 11.333333
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H38").Value = 375.33334
$ws.Range("I38").Value = 375.33334
$ws.Range("K38").Value = 1126.00002
$ws.Range("M38").Value = -754.0000199999999
$ws.Range("H98").Value = 2850.182
$ws.Range("I98").Value = 2390.2
$ws.Range("J98").Value = 7450
$ws.Range("K98").Value = 2390.2
$ws.Range("L98").Value = 7450
$ws.Range("M98").Value = -892.1999999999998
$ws.Range("N98").Value = -10446
$ws.Range("H122").Value = 2850.182
$ws.Range("I122").Value = 2390.2
$ws.Range("J122").Value = 7450
$ws.Range("K122").Value = 7170.599999999999
$ws.Range("L122").Value = 22350
$ws.Range("M122").Value = -4720.599999999999
$ws.Range("N122").Value = -27250
$ws.Range("H137").Value = 2043
$ws.Range("J137").Value = 2937.4
$ws.Range("L137").Value = 8812.200000000001
$ws.Range("N137").Value = -13912.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1694.3636
$ws.Range("I61").Value = 1633.8
$ws.Range("K61").Value = 1633.8
$ws.Range("M61").Value = -1421.8
$ws.Range("H63").Value = 7552
$ws.Range("I63").Value = 104
$ws.Range("K63").Value = 104
$ws.Range("M63").Value = 582
$ws.Range("H66").Value = 7552
$ws.Range("I66").Value = 104
$ws.Range("K66").Value = 520
$ws.Range("M66").Value = 2912
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H101").Value = 30602
$ws.Range("J101").Value = 30602
$ws.Range("L101").Value = 30602
$ws.Range("N101").Value = -37092
$ws.Range("H113").Value = 89666.664
$ws.Range("J113").Value = 89666.664
$ws.Range("L113").Value = 89666.664
$ws.Range("N113").Value = -98344.664
$ws.Range("H132").Value = 5990.6
$ws.Range("I132").Value = 4580.647
$ws.Range("J132").Value = 8986.75
$ws.Range("K132").Value = 13741.941
$ws.Range("L132").Value = 26960.25
$ws.Range("M132").Value = -11211.941
$ws.Range("N132").Value = -32020.25
$ws.Range("H136").Value = 1694.3636
$ws.Range("I136").Value = 1633.8
$ws.Range("K136").Value = 4901.4
$ws.Range("M136").Value = -2351.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 750
$ws.Range("I22").Value = 750
$ws.Range("K22").Value = 750
$ws.Range("M22").Value = -400
$ws.Range("H28").Value = 17999
$ws.Range("J28").Value = 17999
$ws.Range("L28").Value = 17999
$ws.Range("N28").Value = -18489
$ws.Range("H57").Value = 8000
$ws.Range("J57").Value = 8000
$ws.Range("L57").Value = 8000
$ws.Range("N57").Value = -9120
$ws.Range("H105").Value = 1687.1111
$ws.Range("I105").Value = 1199.6666
$ws.Range("K105").Value = 1199.6666
$ws.Range("M105").Value = 547.3334
$ws.Range("H134").Value = 1825.5834
$ws.Range("I134").Value = 1460.9
$ws.Range("K134").Value = 4382.700000000001
$ws.Range("M134").Value = -1847.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2368501.2
$ws.Range("I4").Value = 789501.4399999999
$ws.Range("K4").Value = 2368504.32
$ws.Range("M4").Value = -2368392.32
$ws.Range("H11").Value = 404.33334
$ws.Range("I11").Value = 284.8
$ws.Range("K11").Value = 854.4000000000001
$ws.Range("M11").Value = -714.4000000000001
$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 200
$ws.Range("K16").Value = 600
$ws.Range("M16").Value = -427
$ws.Range("H34").Value = 5274.923
$ws.Range("J34").Value = 5274.923
$ws.Range("L34").Value = 15824.769
$ws.Range("N34").Value = -15992.769
$ws.Range("H39").Value = 6575.5
$ws.Range("J39").Value = 6472.778
$ws.Range("L39").Value = 19418.334
$ws.Range("N39").Value = -20006.334
$ws.Range("H55").Value = 2566.5
$ws.Range("J55").Value = 4800
$ws.Range("L55").Value = 14400
$ws.Range("N55").Value = -14754
$ws.Range("H68").Value = 949.5
$ws.Range("J68").Value = 949.5
$ws.Range("L68").Value = 2848.5
$ws.Range("N68").Value = -4470.5
$ws.Range("H71").Value = 949.5
$ws.Range("J71").Value = 949.5
$ws.Range("L71").Value = 8545.5
$ws.Range("N71").Value = -16657.5
$ws.Range("H99").Value = 2761.625
$ws.Range("I99").Value = 1100
$ws.Range("K99").Value = 3300
$ws.Range("M99").Value = -1054
$ws.Range("H132").Value = 650
$ws.Range("I132").Value = 300
$ws.Range("K132").Value = 2700
$ws.Range("M132").Value = -170

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5813
$ws.Range("I80").Value = 5341.4
$ws.Range("J80").Value = 6599
$ws.Range("K80").Value = 5341.4
$ws.Range("L80").Value = 6599
$ws.Range("M80").Value = -4343.4
$ws.Range("N80").Value = -8595
$ws.Range("H83").Value = 5813
$ws.Range("I83").Value = 5341.4
$ws.Range("J83").Value = 6599
$ws.Range("K83").Value = 26707
$ws.Range("L83").Value = 32995
$ws.Range("M83").Value = -21715
$ws.Range("N83").Value = -42979
$ws.Range("H122").Value = 1952.2222
$ws.Range("I122").Value = 1952.2222
$ws.Range("K122").Value = 5856.6666
$ws.Range("M122").Value = -3406.6666
$ws.Range("H132").Value = 3304.348
$ws.Range("I132").Value = 3220.1
$ws.Range("K132").Value = 9660.299999999999
$ws.Range("M132").Value = -7130.299999999999
$ws.Range("H140").Value = 78000
$ws.Range("J140").Value = 78000
$ws.Range("L140").Value = 78000
$ws.Range("N140").Value = -88360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2184
$ws.Range("I22").Value = 847.25
$ws.Range("J22").Value = 3966.3333
$ws.Range("K22").Value = 847.25
$ws.Range("L22").Value = 3966.3333
$ws.Range("M22").Value = -552.25
$ws.Range("N22").Value = -4556.3333
$ws.Range("H25").Value = 20004
$ws.Range("J25").Value = 22508
$ws.Range("L25").Value = 22508
$ws.Range("N25").Value = -22968
$ws.Range("H27").Value = 2184
$ws.Range("I27").Value = 847.25
$ws.Range("J27").Value = 3966.3333
$ws.Range("K27").Value = 847.25
$ws.Range("L27").Value = 3966.3333
$ws.Range("M27").Value = -740.25
$ws.Range("N27").Value = -4180.3333
$ws.Range("H46").Value = 1382.5834
$ws.Range("I46").Value = 1104.375
$ws.Range("K46").Value = 1104.375
$ws.Range("M46").Value = -916.375
$ws.Range("H82").Value = 2272
$ws.Range("I82").Value = 659
$ws.Range("J82").Value = 3885
$ws.Range("K82").Value = 659
$ws.Range("L82").Value = 3885
$ws.Range("M82").Value = -298
$ws.Range("N82").Value = -4607
$ws.Range("H85").Value = 2272
$ws.Range("I85").Value = 659
$ws.Range("J85").Value = 3885
$ws.Range("K85").Value = 659
$ws.Range("L85").Value = 3885
$ws.Range("M85").Value = 589
$ws.Range("N85").Value = -6381
$ws.Range("H94").Value = 84999.5
$ws.Range("J94").Value = 84999.5
$ws.Range("L94").Value = 84999.5
$ws.Range("N94").Value = -86351.5
$ws.Range("H101").Value = 6965
$ws.Range("J101").Value = 6965
$ws.Range("L101").Value = 6965
$ws.Range("N101").Value = -13455
$ws.Range("H104").Value = 16185
$ws.Range("J104").Value = 16185
$ws.Range("L104").Value = 16185
$ws.Range("N104").Value = -23173
$ws.Range("H132").Value = 4153.8887
$ws.Range("I132").Value = 3231
$ws.Range("K132").Value = 9693
$ws.Range("M132").Value = -7163

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value = 239
$ws.Range("I100").Value = 239
$ws.Range("K100").Value = 478
$ws.Range("M100").Value = 63
$ws.Range("H103").Value = 48333
$ws.Range("J103").Value = 48333
$ws.Range("L103").Value = 48333
$ws.Range("N103").Value = -50677
$ws.Range("H117").Value = 61667
$ws.Range("J117").Value = 61667
$ws.Range("L117").Value = 61667
$ws.Range("N117").Value = -70845
$ws.Range("H122").Value = 1298
$ws.Range("I122").Value = 1298
$ws.Range("K122").Value = 3894
$ws.Range("M122").Value = -1444
